# Append ".global" to the existing header texts in C1:H1 before shifting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = $ws.Range("C1").Value() + ".global"
$ws.Range("D1").Value = $ws.Range("D1").Value() + ".global"
$ws.Range("E1").Value = $ws.Range("E1").Value() + ".global"
$ws.Range("F1").Value = $ws.Range("F1").Value() + ".global"
$ws.Range("G1").Value = $ws.Range("G1").Value() + ".global"
$ws.Range("H1").Value = $ws.Range("H1").Value() + ".global"

# Delete column B ("Year of Treatment"), shifting C:H left to B:G
$ws.Columns("B").Delete()
